$wb = $excel.ActiveWorkbook

# 1) Rename sheet "Disenfranchisement rate" -> "Disenfranchisement Rate"
$rateSheet = $wb.Worksheets.Item("Disenfranchisement rate")
$rateSheet.Name = "Disenfranchisement Rate"

# 2) Fix capitalization of the label in the International sheet
$intlSheet = $wb.Worksheets.Item("International")
$intlSheet.Range("D5").Value = "Individuals Incarcerated *"

# 3) Make "Disenfranchisement Rate" the active/selected sheet/tab
$rateSheet.Activate()
